# Update temperature (D) and uncertainty (E) columns with recomputed
# bootstrap values (feat: histograms of temperature distribution).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = 10.1403762618524;  E = 0.3401649376635241 },
    @{ Row = 3;  D = 10.0144558299543;  E = 0.3295199693804857 },
    @{ Row = 4;  D = 10.97473529764362; E = 0.3392952363579432 },
    @{ Row = 5;  D = 11.27740638296911; E = 0.352369311859356 },
    @{ Row = 6;  D = 11.61404378698954; E = 0.3363480743703146 },
    @{ Row = 7;  D = 12.42598771940949; E = 0.3776266911159403 },
    @{ Row = 8;  D = 12.25033722608613; E = 0.3298622200525696 },
    @{ Row = 9;  D = 13.28187451693018; E = 0.4918318573854879 },
    @{ Row = 10; D = 13.04722761421902; E = 0.3532769108607338 },
    @{ Row = 11; D = 14.3623631510613;  E = 0.4785707621321235 },
    @{ Row = 12; D = 14.02456987818324; E = 0.3981500633966313 },
    @{ Row = 13; D = 15.31754357659524; E = 0.4992000233112369 },
    @{ Row = 14; D = 14.79142891025522; E = 0.4317532294218546 },
    @{ Row = 15; D = 15.96220731988605; E = 0.4902791803780185 },
    @{ Row = 16; D = 15.49040288444248; E = 0.4565458398706101 },
    @{ Row = 17; D = 16.52404572044733; E = 0.4898259902740931 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
